# ============================================================================
# [ADDITIONAL SCRAPING] added scraping code for extra bowling attributes and
# excel sheets.
#
# 1. Clean up left-over blank placeholder cells on "ODI Batting Extra" that
#    the (now extended) scraper no longer emits.
# 2. Add the new "ODI Bowling Extra" sheet (scraped MAIDEN_OVERS /
#    PERCENT_WICKETS_OF_ALL per match) at the end of the workbook.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "ODI Batting Extra" - drop the empty inlineStr placeholder cells
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("ODI Batting Extra")

$ws4.Range("C2:E2").ClearContents()
$ws4.Range("B4:E4").ClearContents()
$ws4.Range("C5:E5").ClearContents()
$ws4.Range("C6:E6").ClearContents()
$ws4.Range("B9:E9").ClearContents()
$ws4.Range("E11").ClearContents()
$ws4.Range("B13:E13").ClearContents()
$ws4.Range("E14").ClearContents()
$ws4.Range("B16:E16").ClearContents()
$ws4.Range("C17:E17").ClearContents()
$ws4.Range("B19:E19").ClearContents()
$ws4.Range("C20:E20").ClearContents()
$ws4.Range("E21").ClearContents()

# ---------------------------------------------------------------------------
# 2. Add "ODI Bowling Extra" as the new last sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = "ODI Bowling Extra"

# Every value scraped onto this sheet is text (even the numeric-looking
# MATCH_CODE / MAIDEN_OVERS / percentage columns), so force text formatting
# before writing any values.
$ws5.Range("A1:C21").NumberFormat = "@"

$header = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($col = 1; $col -le 3; $col++) {
    $ws5.Cells.Item(1, $col).Value = $header[$col - 1]
}

$data = @(
    @("4119", "0", "20.00%"),
    @("4122", "", ""),
    @("4182", "0", "10.00%"),
    @("4210", "0", ""),
    @("4215", "0", ""),
    @("4231", "0", ""),
    @("4232", "0", ""),
    @("4233", "", ""),
    @("4261", "0", ""),
    @("4376", "0", ""),
    @("4413", "0", ""),
    @("4414", "", ""),
    @("4417", "0", ""),
    @("4449", "0", ""),
    @("4450", "", ""),
    @("4451", "2", ""),
    @("4463", "0", "10.00%"),
    @("4464", "", ""),
    @("4480", "0", "10.00%"),
    @("4482", "0", "10.00%")
)

$row = 2
foreach ($rec in $data) {
    $ws5.Cells.Item($row, 1).Value = $rec[0]
    $ws5.Cells.Item($row, 2).Value = $rec[1]
    $ws5.Cells.Item($row, 3).Value = $rec[2]
    $row++
}

# Match the header styling (bold, centered, bordered) already used by the
# other sheets' header rows, by copying the format from an existing header.
$wb.Worksheets.Item("ODI Batting Extra").Range("A1").Copy()
$ws5.Range("A1:C1").PasteSpecial(-4122)
